# Append a new data row (row 96) to each of the four sheets, mirroring the
# structure of the existing rows (time, length, ID, actual length, checksum,
# and their decimal equivalents).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Add-Row96($ws, $timeValue, $bText, $cText, $dText, $eText, $fValue, $gValue, $hValue, $iValue) {
    $row = 96

    $ws.Cells.Item($row, 1).Value = $timeValue
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 2).Value = $bText
    $ws.Cells.Item($row, 3).Value = $cText
    $ws.Cells.Item($row, 4).Value = $dText
    $ws.Cells.Item($row, 5).Value = $eText

    $ws.Cells.Item($row, 6).Value = $fValue
    $ws.Cells.Item($row, 7).Value = $gValue
    $ws.Cells.Item($row, 8).Value = $hValue
    $ws.Cells.Item($row, 9).Value = $iValue
}

# Sheet 1 : FE_LFT_#1
$ws1 = $wb.Worksheets.Item(1)
Add-Row96 $ws1 45882.49700231481 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x04" "0xf" 380 759863127514710945038336.0 260 15

# Sheet 2 : FE_LFT_#2
$ws2 = $wb.Worksheets.Item(2)
Add-Row96 $ws2 45882.49700231481 "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x14" "0xe" 400 568432987514711010443264.0 276 14

# Sheet 3 : FE_PLT_#1
$ws3 = $wb.Worksheets.Item(3)
Add-Row96 $ws3 45882.49700231481 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x5C" "0x3" 110 568631262647113970876416.0 92 3

# Sheet 4 : FE_PLT_#2
$ws4 = $wb.Worksheets.Item(4)
Add-Row96 $ws4 45882.49700231481 "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x5B" "0x3" 110 985046333984776009023488.0 91 3
